$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 99 (pushes the existing row 99 -> 100, and 100 -> 101)
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new weekly record (same
# market/category metadata as the surrounding rows, new date & price figures)
$ws.Cells.Item(99, 1).Value = 9
$ws.Cells.Item(99, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(99, 3).Value = "Metropolitana"
$ws.Cells.Item(99, 4).Value = 44746
$ws.Cells.Item(99, 5).Value = 13
$ws.Cells.Item(99, 6).Value = 100114007
$ws.Cells.Item(99, 7).Value = "Jengibre"
$ws.Cells.Item(99, 8).Value = "Sin especificar"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 700
$ws.Cells.Item(99, 11).Value = 14000
$ws.Cells.Item(99, 12).Value = 15000
$ws.Cells.Item(99, 13).Value = 14500
$ws.Cells.Item(99, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(99, 15).Value = "Perú"
$ws.Cells.Item(99, 16).Value = 1115
$ws.Cells.Item(99, 17).Value = 13
$ws.Cells.Item(99, 18).Value = "Hortaliza"
